# Insert a new weekly price record as row 12 on the single worksheet,
# pushing the existing row 12 (and everything below it) down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(12).Insert()

$ws.Range("A12").Value = 7
$ws.Range("B12").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C12").Value = "Ñuble"
$ws.Range("D12").Value = 45114
$ws.Range("E12").Value = 16
$ws.Range("F12").Value = 100112043
$ws.Range("G12").Value = "Pepino dulce"
$ws.Range("H12").Value = "Cultivar IV Región"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 30
$ws.Range("K12").Value = 15000
$ws.Range("L12").Value = 15000
$ws.Range("M12").Value = 15000
$ws.Range("N12").Value = "$/bandeja 18 kilos"
$ws.Range("O12").Value = "Provincia de Limarí"
$ws.Range("P12").Value = 833
$ws.Range("Q12").Value = 18
$ws.Range("R12").Value = "Hortaliza"
